# 21 Dec - "This time for sure! Presto!!"
# Updates the Regional Bed Avaliability and Hospital COVID Census sheets
# with the latest daily figures.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Regional Bed Avaliability")
$ws2 = $wb.Worksheets.Item("Hospital COVID Census")

# ---------------------------------------------------------------------------
# Sheet 1: "Regional Bed Avaliability" (A2:G8) - new counts for 12/21
# ---------------------------------------------------------------------------
$regionData = @(
    @{Row=2; B=430;  C=2313; D=0;  E=81;  F=555;  G=0}
    @{Row=3; B=151;  C=1271; D=0;  E=71;  F=218;  G=0}
    @{Row=4; B=129;  C=1194; D=0;  E=78;  F=329;  G=0}
    @{Row=5; B=92;   C=844;  D=0;  E=16;  F=147;  G=0}
    @{Row=6; B=90;   C=872;  D=0;  E=63;  F=277;  G=0}
    @{Row=7; B=158;  C=840;  D=23; E=96;  F=230;  G=27}
    @{Row=8; B=1050; C=7334; D=23; E=405; F=1756; G=27}
)

# Comma Style (0 decimals) - matches the workbook's existing custom numFmt 164
$commaFmt = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

foreach ($d in $regionData) {
    $r = $d.Row

    # Write the new values.
    $ws1.Cells.Item($r, 2).Value = $d.B   # B - Occupied ICU
    $ws1.Cells.Item($r, 3).Value = $d.C   # C - Occupied Medical/Surgical
    $ws1.Cells.Item($r, 4).Value = $d.D   # D - Occupied Alternate Medical Site
    $ws1.Cells.Item($r, 5).Value = $d.E   # E - Available ICU
    $ws1.Cells.Item($r, 6).Value = $d.F   # F - Available Medical/Surgical
    $ws1.Cells.Item($r, 7).Value = $d.G   # G - Available Alternate Medical Site

    # Columns B & C: apply the comma number format (font/weight is already
    # correct for the row).
    $ws1.Range("B$r`:C$r").NumberFormat = $commaFmt

    # Columns D & G: switch to the Arial body font used by the rest of the
    # table (keeps the General number format).
    $ws1.Range("D$r").Font.Name = "Arial"
    $ws1.Range("G$r").Font.Name = "Arial"

    # Columns E & F: Arial font plus the comma number format.
    $ws1.Range("E$r`:F$r").Font.Name = "Arial"
    $ws1.Range("E$r`:F$r").NumberFormat = $commaFmt

    # Row 8 ("Total") carries the bold weight used elsewhere in that row.
    if ($r -eq 8) {
        $ws1.Range("D$r").Font.Bold = $true
        $ws1.Range("G$r").Font.Bold = $true
        $ws1.Range("E$r`:F$r").Font.Bold = $true
    }
}

$ws1.Range("G2:G8").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Hospital COVID Census" (C:D) - updated hospitalization / ICU
# census counts for 12/21
# ---------------------------------------------------------------------------
$hospData = @(
    @{Row=3;  C=21;  D=2}
    @{Row=5;  C=7;   D=2}
    @{Row=6;  C=121; D=17}
    @{Row=7;  C=8;   D=0}
    @{Row=8;  C=10;  D=0}
    @{Row=9;  C=36;  D=3}
    @{Row=10; C=21;  D=5}
    @{Row=11; C=12;  D=3}
    @{Row=13; C=72;  D=22}
    @{Row=14; C=31;  D=3}
    @{Row=15; C=4;   D=2}
    @{Row=16; C=78;  D=17}
    @{Row=17; C=32;  D=5}
    @{Row=18; C=92;  D=27}
    @{Row=19; C=28;  D=5}
    @{Row=20; C=42;  D=5}
    @{Row=21; C=24;  D=1}
    @{Row=22; C=34;  D=9}
    @{Row=24; C=11;  D=1}
    @{Row=27; C=2;   D=0}
    @{Row=28; C=12;  D=2}
    @{Row=29; C=39;  D=5}
    @{Row=30; C=14;  D=2}
    @{Row=31; C=38;  D=8}
    @{Row=32; C=17;  D=4}
    @{Row=33; C=46;  D=8}
    @{Row=34; C=15;  D=2}
    @{Row=35; C=65;  D=15}
    @{Row=37; C=34;  D=5}
    @{Row=38; C=51;  D=6}
    @{Row=39; C=15;  D=7}
    @{Row=42; C=116; D=40}
    @{Row=43; C=28;  D=7}
    @{Row=44; C=16;  D=4}
    @{Row=46; C=29;  D=11}
    @{Row=48; C=47;  D=6}
    @{Row=49; C=46;  D=7}
    @{Row=50; C=35;  D=6}
    @{Row=52; C=4;   D=0}
    @{Row=54; C=37;  D=5}
    @{Row=55; C=55;  D=9}
    @{Row=57; C=52;  D=12}
    @{Row=58; C=8;   D=0}
    @{Row=59; C=51;  D=9}
    @{Row=60; C=37;  D=7}
    @{Row=61; C=36;  D=5}
    @{Row=62; C=29;  D=7}
    @{Row=63; C=41;  D=6}
    @{Row=64; C=29;  D=4}
    @{Row=65; C=8;   D=1}
    @{Row=66; C=32;  D=16}
    @{Row=67; C=23;  D=0}
    @{Row=68; C=59;  D=15}
    @{Row=69; C=44;  D=29}
    @{Row=70; C=52;  D=6}
)

foreach ($d in $hospData) {
    $ws2.Cells.Item($d.Row, 3).Value = $d.C
    $ws2.Cells.Item($d.Row, 4).Value = $d.D
}

$ws2.Range("F16").Select()
